$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows for PID 3 and PID 4, mirroring the pattern of PID 1 and PID 2
$newRows = @(
    @(3, "R1", "ND"),
    @(3, "R2", "D"),
    @(3, "R3", "ND"),
    @(3, "R4", "D"),
    @(4, "R1", "D"),
    @(4, "R2", "ND"),
    @(4, "R3", "D"),
    @(4, "R4", "ND")
)

$startRow = 12
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
}

$ws.Range("N15").Select() | Out-Null
